# Add two new rows to the DS-AlgoQns tracker sheet:
#   Row 6: Strings / Valid Palindrome
#   Row 7: Searching And Sorting / First Bad version

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 6: Valid Palindrome ----
$ws.Range("A6").Value = "Strings"
$ws.Range("B6").Value = "Valid Palindrome"

$ws.Range("C6").Value = "https://leetcode.com/explore/interview/card/top-interview-questions-easy/127/strings/883/"
$ws.Hyperlinks.Add($ws.Range("C6"), "https://leetcode.com/explore/interview/card/top-interview-questions-easy/127/strings/883/")

$ws.Range("D6").Value = "Use two pointer approach if characters are equal or not"
$ws.Range("D6").WrapText = $true

$ws.Range("E6").Value = "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoLibrary/Strings/Palindrome.cs"
$ws.Hyperlinks.Add($ws.Range("E6"), "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoLibrary/Strings/Palindrome.cs")

$ws.Range("F6").Value = "Easy - make sure to skip any non alphanumeric characters"
$ws.Range("G6").Value = "O(N)"
$ws.Range("H6").Value = "O(1)"

# restore the Hyperlink cell style (Hyperlinks.Add re-creates a near duplicate style)
$ws.Range("C6").Style = "Hyperlink"
$ws.Range("E6").Style = "Hyperlink"

$ws.Rows.Item(6).RowHeight = 29

# ---- Row 7: First Bad version ----
$ws.Range("A7").Value = "Searching And Sorting"
$ws.Range("B7").Value = "First Bad version"

$ws.Range("C7").Value = "https://leetcode.com/explore/interview/card/top-interview-questions-easy/96/sorting-and-searching/774/"
$ws.Hyperlinks.Add($ws.Range("C7"), "https://leetcode.com/explore/interview/card/top-interview-questions-easy/96/sorting-and-searching/774/")

$ws.Range("D7").Value = "Use Binary search approach to find bad version until both pointers resolve to 1"
$ws.Range("D7").WrapText = $true

$ws.Range("G7").Value = "O(logN)"
$ws.Range("F7").Value = "Easy - Binary Search"
$ws.Range("H7").Value = "O(1)"

$ws.Range("C7").Style = "Hyperlink"

$ws.Rows.Item(7).RowHeight = 29

# Match the saved selection state from the authored workbook
[void]$ws.Range("E7").Select()
